# Realestate Update resale numbers 2023-06-07 20:39
# Appends a new data row (row 28) to the CityResaleNum sheet with the
# latest resale-number snapshot, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# Columns A-D are plain text in this sheet (date/time/weekday/week-as-text).
# Force text formatting before assignment so Excel doesn't reinterpret
# "2023-06-07" as a date serial or "23" as a number, then clear the
# number-format override (ClearFormats) so the cell keeps its text value
# but reverts to the sheet's default (unstyled) cell format, matching the
# rest of the data rows.
function Set-TextCell($r, $c, $value) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell $row 1 "2023-06-07"
Set-TextCell $row 2 "20:36:44"
Set-TextCell $row 3 "Wednesday"
Set-TextCell $row 4 "23"

# Columns E-T are numeric city resale-number values.
$ws.Cells.Item($row, 5).Value  = 116944
$ws.Cells.Item($row, 6).Value  = 134160
$ws.Cells.Item($row, 7).Value  = 160143
$ws.Cells.Item($row, 8).Value  = 130589
$ws.Cells.Item($row, 9).Value  = 175471
$ws.Cells.Item($row, 10).Value = 112816
$ws.Cells.Item($row, 11).Value = 200909
$ws.Cells.Item($row, 12).Value = 220950
$ws.Cells.Item($row, 13).Value = 172482
$ws.Cells.Item($row, 14).Value = 119837
$ws.Cells.Item($row, 15).Value = 38576
$ws.Cells.Item($row, 16).Value = 34501
$ws.Cells.Item($row, 17).Value = 50714
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36936
$ws.Cells.Item($row, 20).Value = -1
